# Applies the AudioBaseline.xlsx edit described by the commit:
#   "SeparateData: added alternative self-report analysis, added plots
#    PlotResults: plot offer as function of emotion strength"
#
# Sheet1 holds a long (row id, participant/clip id, emotion, score) table
# that ends in a sentinel row of "NULL" strings (previously row 156). The
# edit appends 5 more participant blocks (5 emotions each = 25 rows) of
# self-report data before that trailing sentinel row, which pushes the
# sentinel row from 156 down to 181.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 25 rows at row 156 (inherits formatting/style from the row above,
# same as Excel native row insertion), pushing the old NULL sentinel row
# down from row 156 to row 181.
$ws.Rows("156:180").Insert()

# New data rows: (row id, participant/clip id, emotion, score)
$newRows = @(
    @(156, 1032, "Neutral", 0.1823),
    @(157, 1032, "Happy", 0),
    @(158, 1032, "Sad", 0.0001),
    @(159, 1032, "Angry", 0.0134),
    @(160, 1032, "Fear", 0.8041),
    @(161, 1033, "Neutral", 0.1769),
    @(162, 1033, "Happy", 0.0007),
    @(163, 1033, "Sad", 0.7921),
    @(164, 1033, "Angry", 0.0001),
    @(165, 1033, "Fear", 0.0302),
    @(166, 1035, "Neutral", 0.2153),
    @(167, 1035, "Happy", 0),
    @(168, 1035, "Sad", 0.7757),
    @(169, 1035, "Angry", 0),
    @(170, 1035, "Fear", 0.009),
    @(171, 1039, "Neutral", 0.942),
    @(172, 1039, "Happy", 0),
    @(173, 1039, "Sad", 0.0512),
    @(174, 1039, "Angry", 0),
    @(175, 1039, "Fear", 0.0069),
    @(176, 1038, "Neutral", 0.0241),
    @(177, 1038, "Happy", 0.0163),
    @(178, 1038, "Sad", 0.0014),
    @(179, 1038, "Angry", 0.0485),
    @(180, 1038, "Fear", 0.9097),
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update the view: it was scrolled down to the newly-added tail of data,
# with J176 the active cell of the selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 169
$win.ScrollColumn = 1
$ws.Range("J176").Select()
